# Apply "Point of overcoming" noise-statistics update to the "withoutNoise" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("withoutNoise")

# --- Row 3 (Permanent / Point of overcoming) : all four stat columns updated ---
$ws.Range("B3").Value = 328
$ws.Range("C3").Value = 328
$ws.Range("G3").Value = 350
$ws.Range("H3").Value = 350
$ws.Range("L3").Value = 329
$ws.Range("M3").Value = 329
$ws.Range("Q3").Value = 327
$ws.Range("R3").Value = 327

# --- Row 13 (Temporary / Point of overcoming) ---
$ws.Range("B13").Value = 330
$ws.Range("C13").Value = 330
$ws.Range("Q13").Value = 327
$ws.Range("R13").Value = 326

# --- Row 23 (Shifted / Point of overcoming) ---
$ws.Range("B23").Value = 334
$ws.Range("C23").Value = 333
$ws.Range("Q23").Value = 332
$ws.Range("R23").Value = 331

# --- Row 33 (Outlier / Point of overcoming) ---
$ws.Range("B33").Value = 306
$ws.Range("C33").Value = 306
$ws.Range("L33").Value = 306
$ws.Range("M33").Value = 306
$ws.Range("Q33").Value = 305
$ws.Range("R33").Value = 305

# --- The "Col" (G:H) and "Sym" (L:M) blocks for the Temporary, Shifted and
#     Outlier groups no longer mirror the Row/Diag statistics now that the
#     point of overcoming uses modelling-noise statistics, so those cells
#     (point-of-overcoming row plus the four X[Q+n] rows beneath it) are
#     cleared entirely.
$ws.Range("G13:G17").ClearContents()
$ws.Range("H13:H17").ClearContents()
$ws.Range("L13:L17").ClearContents()
$ws.Range("M13:M17").ClearContents()

$ws.Range("G23:G27").ClearContents()
$ws.Range("H23:H27").ClearContents()
$ws.Range("L23:L27").ClearContents()
$ws.Range("M23:M27").ClearContents()

$ws.Range("G33:G37").ClearContents()
$ws.Range("H33:H37").ClearContents()
